# "Week 25" attendance data was still blank (npc's / dialogue work week) -
# fill in the GL-uren hours for Dees, Thomas and Arne, same pattern as the
# other weeks: Tue=4, Wed=2, Thur=4, Fri=8.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 25")
$ws.Activate()

$ws.Range("C3:E3").Value = 4
$ws.Range("C4:E4").Value = 2
$ws.Range("C5:E5").Value = 4
$ws.Range("C6:E6").Value = 8

$ws.Range("G8").Select()
